$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: drop the (redundant) explicit style from A2/B2 ------------------
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Style = "Normal"

# --- Row 26: two new cells (H26/I26) ----------------------------------------
$ws.Range("H26").Value = 7669
$ws.Range("I26").Formula = "=7669/3600"

# --- Row 27: drop the old "Total" label (moved to row 28) and the old D27 ---
#     helper column; repoint the E27 time-conversion formula at C27 directly.
$ws.Range("B27").Clear()
$ws.Range("D27").Clear()
$ws.Range("E27").Formula = "=C27/86400"

# --- Row 28: new "Total" row -------------------------------------------------
$ws.Range("C28").Style = "Normal"
$ws.Range("C28").Formula = "=SUM(C2:C27)"

$ws.Range("B28").Value = "Total"
$ws.Range("B28").Font.Bold = $true

$ws.Range("E28").Formula = "=C28/86400"
$ws.Range("E28").NumberFormat = "[h]:mm:ss;@"
$ws.Range("E28").Font.Bold = $true

# --- View state: zoom + new selection ---------------------------------------
$excel.ActiveWindow.Zoom = 130
$null = $ws.Range("H26").Select()
